$d = $word.ActiveDocument

# 1) Merge the split "Sometimes a dataset..." run back together and drop the
#    _GoBack bookmark that was sitting between the two runs. Find matches the
#    paragraph's rendered text (already the post-edit text, since the two
#    runs concatenate to it) and replaces the whole matched range - which
#    spans both runs and the bookmark - with a single run.
$d.Content.Find.Execute(
    "Sometimes a dataset is acquired before coming up w/ a question + we investigate the data to think of one",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sometimes a dataset is acquired before coming up w/ a question + we investigate the data to think of one",
    2) | Out-Null

# 2) Insert the new paragraphs (Data Wrangling section) right after that
#    paragraph by splicing raw OOXML into a collapsed point range positioned
#    just before the paragraph mark of the last paragraph.
$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)

$newParagraphsXml = @'
<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:rPr><w:i/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:u w:val="single"/></w:rPr><w:t>Data Wrangling</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:rPr><w:i/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:tabs><w:tab w:val="clear" w:pos="1530"/><w:tab w:val="num" w:pos="360"/></w:tabs><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve">Can get data from a CSV, a database, an API, scraping a webpage, and </w:t></w:r><w:r><w:t>then</w:t></w:r><w:r><w:t xml:space="preserve"> combine data from </w:t></w:r><w:r><w:t>various sources</w:t></w:r><w:r><w:t xml:space="preserve"> that are in </w:t></w:r><w:r><w:t>different</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>formats</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:tabs><w:tab w:val="clear" w:pos="1530"/><w:tab w:val="num" w:pos="360"/></w:tabs><w:ind w:left="360"/></w:pPr><w:r><w:t>CSV files have no formulas and are much easier to process w/ code, unlike Excel spreadsheets</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:tabs><w:tab w:val="clear" w:pos="1530"/><w:tab w:val="num" w:pos="360"/></w:tabs><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve">Can represent a CSV in python where </w:t></w:r><w:r><w:t>each</w:t></w:r><w:r><w:t xml:space="preserve"> row is</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:tabs><w:tab w:val="clear" w:pos="1530"/><w:tab w:val="num" w:pos="720"/></w:tabs><w:ind w:left="720"/></w:pPr><w:r><w:t>A list, so the entire dataset is a list of lists</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:tabs><w:tab w:val="clear" w:pos="1530"/><w:tab w:val="num" w:pos="720"/></w:tabs><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">A </w:t></w:r><w:r><w:t>dictionary</w:t></w:r><w:r><w:t xml:space="preserve">, which works well if the CSV has a header </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> keys = col names, values = fields</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:tabs><w:tab w:val="clear" w:pos="1530"/><w:tab w:val="num" w:pos="360"/></w:tabs><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve">Can use </w:t></w:r><w:r><w:t>Python’s</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">unicodecsv </w:t></w:r><w:r><w:t>library to load in CSV files</w:t></w:r><w:r><w:t>, or use open to open a file connection and read it in</w:t></w:r></w:p>
'@

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newParagraphsXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($packageXml)

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
